$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics to account for the new trade #43
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.09   # Current Capital
$summary.Range("B4").Value = 0.08      # Total P&L $
$summary.Range("B6").Value = 43        # Total Trades
$summary.Range("B8").Value = 17        # Losing Trades
$summary.Range("B9").Value = 39.53     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.09     # Capital
$status.Range("D4").Value = 43         # Trades
$status.Range("E4").Value = 0.08       # P&L $
$status.Range("F4").Value = 0.09       # P&L %
$status.Range("G4").Value = 39.53      # Win Rate %

# ---------------------------------------------------------------------------
# Append the new, closed trade #43 to both the "All Trades" and
# "MarketMaking" sheets (they track the same trade log).
# ---------------------------------------------------------------------------
$sheetsToAppend = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetsToAppend) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 44

    $ws.Cells.Item($row, 1).Value = 43

    # Force the date/time-looking strings to remain plain text instead of
    # being auto-converted to Excel date/time serials.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "12:47:49"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.11
    $ws.Cells.Item($row, 7).Value = 0.1
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -9.0909
    $ws.Cells.Item($row, 10).Value = -0.01
    $ws.Cells.Item($row, 11).Value = 100.09
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
